$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and 1h volume/change (E) columns for each coin row.
# Column D values that parse as a plain number (e.g. "96.29") are written
# with a leading apostrophe so Excel stores them as text, matching the
# workbook's existing text-formatted "Price" column (values like "43.664.14"
# already cannot be numbers because of the thousands-dot formatting).

$ws.Cells.Item(2, "D").Value = '43.571.05'
$ws.Cells.Item(2, "E").Value = '  -0.54%  '
$ws.Cells.Item(3, "D").Value = '2.288.18'
$ws.Cells.Item(3, "E").Value = '  -1.24%  '
$ws.Cells.Item(4, "E").Value = '  +0.08%  '
$ws.Cells.Item(5, "D").Value = '''96.29'
$ws.Cells.Item(5, "E").Value = '  +2.47%  '
$ws.Cells.Item(6, "D").Value = '''268.35'
$ws.Cells.Item(6, "E").Value = '  -0.81%  '
$ws.Cells.Item(7, "D").Value = '''0.615'
$ws.Cells.Item(7, "E").Value = '  -2.20%  '
$ws.Cells.Item(8, "E").Value = '  -0.05%  '
$ws.Cells.Item(9, "D").Value = '''0.608'
$ws.Cells.Item(9, "E").Value = '  -2.52%  '
$ws.Cells.Item(10, "D").Value = '''45.36'
$ws.Cells.Item(10, "E").Value = '  +0.95%  '
$ws.Cells.Item(11, "E").Value = '  +0.08%  '
$ws.Cells.Item(12, "D").Value = '''7.87'
$ws.Cells.Item(12, "E").Value = '  -3.74%  '
$ws.Cells.Item(13, "E").Value = '  +0.92%  '
$ws.Cells.Item(14, "D").Value = '2.631.55'
$ws.Cells.Item(14, "E").Value = '  -1.11%  '
$ws.Cells.Item(15, "D").Value = '''15.36'
$ws.Cells.Item(15, "E").Value = '  -0.47%  '
$ws.Cells.Item(16, "D").Value = '''0.848'
$ws.Cells.Item(16, "E").Value = '  -1.99%  '
$ws.Cells.Item(17, "D").Value = '2.288.39'
$ws.Cells.Item(17, "E").Value = '  -0.85%  '
$ws.Cells.Item(18, "D").Value = '43.553.54'
$ws.Cells.Item(18, "E").Value = '  -0.60%  '
$ws.Cells.Item(19, "E").Value = '  +2.46%  '
$ws.Cells.Item(20, "D").Value = '''6.19'
$ws.Cells.Item(20, "E").Value = '  -1.83%  '
$ws.Cells.Item(21, "D").Value = '''72.05'
$ws.Cells.Item(21, "E").Value = '  +0.58%  '
$ws.Cells.Item(22, "D").Value = '''2.59'
$ws.Cells.Item(22, "E").Value = '  +12.67%  '
$ws.Cells.Item(23, "D").Value = '''232.41'
$ws.Cells.Item(23, "E").Value = '  -1.92%  '
$ws.Cells.Item(24, "D").Value = '''9.15'
$ws.Cells.Item(24, "E").Value = '  -5.29%  '
$ws.Cells.Item(25, "D").Value = '''2.61'
$ws.Cells.Item(25, "E").Value = '  +4.03%  '
$ws.Cells.Item(26, "E").Value = '  -0.05%  '
$ws.Cells.Item(27, "D").Value = '''11.23'
$ws.Cells.Item(27, "E").Value = '  -1.49%  '
$ws.Cells.Item(28, "E").Value = '  +2.31%  '
$ws.Cells.Item(29, "D").Value = '''39.87'
$ws.Cells.Item(29, "E").Value = '  +1.98%  '
$ws.Cells.Item(30, "D").Value = '''2.23'
$ws.Cells.Item(30, "E").Value = '  -5.80%  '
$ws.Cells.Item(31, "D").Value = '''174.72'
$ws.Cells.Item(31, "E").Value = '  +1.48%  '
$ws.Cells.Item(32, "D").Value = '''21.78'
$ws.Cells.Item(32, "E").Value = '  -3.25%  '
$ws.Cells.Item(33, "D").Value = '''0.0892'
$ws.Cells.Item(33, "E").Value = '  -0.72%  '
$ws.Cells.Item(34, "D").Value = '''5.37'
$ws.Cells.Item(34, "E").Value = '  -3.36%  '
$ws.Cells.Item(35, "D").Value = '''0.126'
$ws.Cells.Item(35, "E").Value = '  -0.74%  '
$ws.Cells.Item(36, "E").Value = '  -4.16%  '
$ws.Cells.Item(37, "D").Value = '''0.0352'
$ws.Cells.Item(37, "E").Value = '  -2.73%  '
$ws.Cells.Item(38, "E").Value = '  -3.67%  '
$ws.Cells.Item(39, "D").Value = '''3.35'
$ws.Cells.Item(39, "E").Value = '  -4.70%  '
$ws.Cells.Item(40, "E").Value = '  +1.68%  '
$ws.Cells.Item(41, "E").Value = '  -0.36%  '
$ws.Cells.Item(42, "D").Value = '''12.30'
$ws.Cells.Item(42, "E").Value = '  -0.04%  '
$ws.Cells.Item(43, "D").Value = '''65.50'
$ws.Cells.Item(43, "E").Value = '  +5.39%  '
$ws.Cells.Item(44, "E").Value = '  +0.37%  '
$ws.Cells.Item(45, "D").Value = '''8.79'
$ws.Cells.Item(45, "E").Value = '  -2.48%  '
$ws.Cells.Item(46, "E").Value = '  -0.94%  '
$ws.Cells.Item(47, "D").Value = '''5.14'
$ws.Cells.Item(47, "E").Value = '  -6.28%  '
$ws.Cells.Item(48, "D").Value = '''96.91'
$ws.Cells.Item(48, "E").Value = '  -3.10%  '
$ws.Cells.Item(49, "E").Value = '  -1.73%  '

# Rows 50 and 51 swap places (WOONetwork <-> TheGraph) with refreshed values.
$ws.Cells.Item(50, "B").Value = 'TheGraph'
$ws.Cells.Item(50, "C").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(50, "D").Value = '''0.185'
$ws.Cells.Item(50, "E").Value = '  +5.96%  '
$ws.Cells.Item(51, "B").Value = 'WOONetwork'
$ws.Cells.Item(51, "C").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Cells.Item(51, "D").Value = '''0.430'
$ws.Cells.Item(51, "E").Value = '  +0.09%  '
